$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Segal, Paul"
$ws.Range("A2").Value = "Gordon, Victorine"
$ws.Range("A3").Value = "Levenson, David"
$ws.Range("A4").Value = "Matzkin, Harold"
$ws.Range("A5").Value = "Shnayder, Leon"
$ws.Range("A6").Value = "Sokolinski, Ilia"
$ws.Range("A7").Value = "Berkovits, Shimshon"

$ws.Range("A2").Select()
